$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column Y ("capacity" is X, "encoded_days" was Y and
# becomes Z). Everything from Y..AI shifts right to Z..AJ automatically.
$ws.Columns("Y").Insert()

# New header for the inserted column.
$ws.Range("Y1").Value = "encoded_year_taken"

# Fill the new column with the encoded year value for every data row (2-62).
# Store it as text (matching the "Normal Year Taken" / S column convention of
# storing these small codes as strings rather than numbers).
$dataRange = $ws.Range("Y2:Y62")
$dataRange.NumberFormat = "@"
$dataRange.Value = "2"
